# Applies the "finished 10 novice users data" commit:
#  - fills in the last novice-user template (Sheet8) with trial data
#  - adds Sheet9 and Sheet10 (novice users #9 and #10), each filled with data
#  - adds Expert 1 (empty data template, not yet filled) and Expert 2 (brand
#    new blank sheet) as placeholders for the not-yet-started expert trials
#  - moves the active/selected tab back to Sheet1

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Finish Sheet8 (was an empty trial template, now gets real measurements)
# ---------------------------------------------------------------------------
$sheet8 = $wb.Worksheets.Item("Sheet8")

$sheet8.Range("C3").Value = 7.6
$sheet8.Range("D3").Value = 25.9
$sheet8.Range("E3").Value = 41.7
$sheet8.Range("F3").Value = 190.8

$sheet8.Range("C4").Value = 9.6
$sheet8.Range("D4").Value = 20.4
$sheet8.Range("E4").Value = 48.4
$sheet8.Range("F4").Value = 199.4

$sheet8.Range("C5").Value = 10
$sheet8.Range("D5").Value = 20.3
$sheet8.Range("E5").Value = 48.7
$sheet8.Range("F5").Value = 202.2

$sheet8.Range("F15").Select()

# ---------------------------------------------------------------------------
# 2. Sheet9 - novice user #9, filled in with the standard template + data
# ---------------------------------------------------------------------------
$sheet9 = $wb.Worksheets.Add($null, $sheet8)
$sheet9.Name = "Sheet9"

$sheet9.Range("B1").Value = "Measurement "
$sheet9.Range("C1").Value = 10
$sheet9.Range("D1").Value = 20
$sheet9.Range("E1").Value = 50
$sheet9.Range("F1").Value = 200
$sheet9.Range("J1").Value = "weight in mg"

$sheet9.Range("A2").Value = "Trial"

$sheet9.Range("A3").Value = 1
$sheet9.Range("C3").Value = 10.9
$sheet9.Range("D3").Value = 15.9
$sheet9.Range("E3").Value = 40.3
$sheet9.Range("F3").Value = 187

$sheet9.Range("A4").Value = 2
$sheet9.Range("C4").Value = 10.2
$sheet9.Range("D4").Value = 20.4
$sheet9.Range("E4").Value = 47.1
$sheet9.Range("F4").Value = 182.6

$sheet9.Range("A5").Value = 3
$sheet9.Range("C5").Value = 11.1
$sheet9.Range("D5").Value = 19.2
$sheet9.Range("E5").Value = 47.1
$sheet9.Range("F5").Value = 186

$sheet9.Range("K27").Select()

# ---------------------------------------------------------------------------
# 3. Sheet10 - novice user #10, completing "10 novice users" (no J1 label
#    on this one, matching the source data exactly)
# ---------------------------------------------------------------------------
$sheet10 = $wb.Worksheets.Add($null, $sheet9)
$sheet10.Name = "Sheet10"

$sheet10.Range("B1").Value = "Measurement "
$sheet10.Range("C1").Value = 10
$sheet10.Range("D1").Value = 20
$sheet10.Range("E1").Value = 50
$sheet10.Range("F1").Value = 200

$sheet10.Range("A2").Value = "Trial"

$sheet10.Range("A3").Value = 1
$sheet10.Range("C3").Value = 8.8
$sheet10.Range("D3").Value = 20.4
$sheet10.Range("E3").Value = 54.1
$sheet10.Range("F3").Value = 193.8

$sheet10.Range("A4").Value = 2
$sheet10.Range("C4").Value = 10.9
$sheet10.Range("D4").Value = 16.5
$sheet10.Range("E4").Value = 41.1
$sheet10.Range("F4").Value = 197.6

$sheet10.Range("A5").Value = 3
$sheet10.Range("C5").Value = 11.8
$sheet10.Range("D5").Value = 16.3
$sheet10.Range("E5").Value = 42
$sheet10.Range("F5").Value = 197.7

$sheet10.Range("F5").Select()

# ---------------------------------------------------------------------------
# 4. Expert 1 - the blank trial template, ready for data entry (not filled
#    in yet as part of this commit)
# ---------------------------------------------------------------------------
$expert1 = $wb.Worksheets.Add($null, $sheet10)
$expert1.Name = "Expert 1"

$expert1.Range("B1").Value = "Measurement "
$expert1.Range("C1").Value = 10
$expert1.Range("D1").Value = 20
$expert1.Range("E1").Value = 50
$expert1.Range("F1").Value = 200
$expert1.Range("J1").Value = "weight in mg"

$expert1.Range("A2").Value = "Trial"

$expert1.Range("A3").Value = 1
$expert1.Range("A4").Value = 2
$expert1.Range("A5").Value = 3

$expert1.Range("C3").Select()

# ---------------------------------------------------------------------------
# 5. Expert 2 - freshly inserted, completely blank so far
# ---------------------------------------------------------------------------
$expert2 = $wb.Worksheets.Add($null, $expert1)
$expert2.Name = "Expert 2"

# ---------------------------------------------------------------------------
# 6. Move the active tab back to Sheet1 (it had drifted to Sheet7)
# ---------------------------------------------------------------------------
$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet1.Activate()
